$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'95.280.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.89%  "

$ws.Range("D3").Value = "'3.591.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.08%  "

$ws.Range("D4").Value = "'2.69"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +31.35%  "

$ws.Range("D5").Value = "'0.999"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").Value = "'221.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.95%  "

$ws.Range("D7").Value = "'635.79"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.19%  "

$ws.Range("D8").Value = "'0.420"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.79%  "

$ws.Range("E9").Value = "  +9.37%  "

$ws.Range("D10").Value = "'0.999"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("D11").Value = "'3.587.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.16%  "

$ws.Range("D12").Value = "'47.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.28%  "

$ws.Range("D13").Value = "'0.213"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.77%  "

$ws.Range("D14").Value = "'0.0000292"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.73%  "

$ws.Range("D15").Value = "'6.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.10%  "

$ws.Range("D16").Value = "'4.260.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.24%  "

$ws.Range("D17").Value = "'95.048.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.97%  "

$ws.Range("D18").Value = "'22.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +22.78%  "

$ws.Range("D19").Value = "'8.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.89%  "

$ws.Range("D20").Value = "'13.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.73%  "

$ws.Range("D21").Value = "'3.583.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.25%  "

$ws.Range("E22").Value = "  +46.66%  "

$ws.Range("E23").Value = "  +2.15%  "

$ws.Range("D24").Value = "'510.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.73%  "

$ws.Range("E25").Value = "  -6.87%  "

$ws.Range("D26").Value = "'125.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +16.72%  "

$ws.Range("D27").Value = "'0.0000201"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.15%  "

$ws.Range("D28").Value = "'6.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.58%  "

$ws.Range("D29").Value = "'3.757.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.87%  "

$ws.Range("D30").Value = "'12.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.19%  "

$ws.Range("E31").Value = "  +1.64%  "

$ws.Range("E32").Value = "  +0.48%  "

$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("D34").Value = "'0.619"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.12%  "

$ws.Range("E35").Value = "  -5.75%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.13%  "

$ws.Range("D37").Value = "'32.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.39%  "

$ws.Range("E38").Value = "  -4.04%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").Value = "'0.530"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.21%  "

$ws.Range("D41").Value = "'7.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.83%  "

$ws.Range("D42").Value = "'8.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.90%  "

$ws.Range("D43").Value = "'579.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.72%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'41.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.94%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0516"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.62%  "

$ws.Range("E46").Value = "  -6.26%  "

$ws.Range("D47").Value = "'0.954"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.93%  "

$ws.Range("D48").Value = "'1.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.95%  "

$ws.Range("D49").Value = "'9.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.69%  "

$ws.Range("D50").Value = "'231.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.78%  "

$ws.Range("E51").Value = "  -0.58%  "
